$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsDeDe = $wb.Worksheets.Item("de-de")

# Status text: "Handed back: in sync with en-US" -> "Ready for handoff"
$wsOverview.Range("E2").Value = "Ready for handoff"
$wsOverview.Range("F2").Value = "Ready for handoff"
$wsZhCn.Range("C2").Value = "Ready for handoff"
$wsDeDe.Range("C2").Value = "Ready for handoff"

# Latest HO Xliff Generate Date / Latest Handback DateTime: 2016-08-25 19:01:21 -> 2016-08-25 19:02:09
$wsOverview.Range("G2").Value = "2016-08-25 19:02:09"
$wsDeDe.Range("H2").Value = "2016-08-25 19:02:09"

# zh-cn Latest Handoff Datetime: 2016-08-25 19:01:16 -> 2016-08-25 19:01:58
$wsZhCn.Range("H2").Value = "2016-08-25 19:01:58"

# Column width changes: 29.9777047293527 -> 17.2159881591797
# (the stored OOXML "width" = ColumnWidth + 5/6, snapped to the nearest 1/6;
#  16.333333333333336 is the ColumnWidth that lands closest on the target)
$wsOverview.Columns("E").ColumnWidth = 16.333333333333336
$wsOverview.Columns("F").ColumnWidth = 16.333333333333336
$wsZhCn.Columns("C").ColumnWidth = 16.333333333333336
$wsDeDe.Columns("C").ColumnWidth = 16.333333333333336
